# Update "想去人数" (want-to-go count) figures pulled from the latest
# generated data (gh-pages output at 456a3b4).
#
# Two worksheets carry the same rows of exhibition data and both need
# their F-column counts refreshed:
#   - "展览"     (Exhibitions)
#   - "全部类型" (All types)

$wb = $excel.ActiveWorkbook

# Map of old value -> new value for column F cells.
$updates = @{
    11863 = 11866
    11786 = 11791
    53    = 56
    5837  = 5841
    3533  = 3534
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 6)  # column F
        $old = $cell.Value2
        if ($updates.ContainsKey($old)) {
            $cell.Value2 = $updates[$old]
        }
    }
}
